$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-11-20 Wednesday" "2024-11-21 Thursday"

Replace-Text "793÷5=" "847÷2="
Replace-Text "418÷7=" "643÷2="
Replace-Text "633÷3=" "157÷8="
Replace-Text "425÷6=" "545÷8="
Replace-Text "140÷7=" "932÷6="
Replace-Text "595÷8=" "982÷3="
Replace-Text "841÷3=" "151÷2="
Replace-Text "762÷4=" "134÷2="
Replace-Text "634÷6=" "926÷4="
Replace-Text "660÷7=" "626÷8="
Replace-Text "843÷3=" "433÷8="
Replace-Text "149÷7=" "312÷5="
Replace-Text "687÷6=" "173÷6="
Replace-Text "576÷4=" "449÷6="
Replace-Text "870÷9=" "408÷9="
Replace-Text "601÷2=" "431÷2="
Replace-Text "460÷7=" "812÷5="
Replace-Text "754÷7=" "215÷5="
Replace-Text "370÷5=" "570÷7="
Replace-Text "816÷8=" "898÷4="
Replace-Text "448÷2=" "163÷5="
Replace-Text "918÷9=" "791÷3="
Replace-Text "502÷7=" "292÷5="
Replace-Text "515÷9=" "554÷2="
Replace-Text "804÷9=" "240÷7="
